$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("T3").Value = 0
$ws.Range("W3").Value = 7.0303
$ws.Range("AB3").Value = 0
$ws.Range("AF3").Value = 0
$ws.Range("AJ3").Value = 17.9394
$ws.Range("B4").Value = 18.174
$ws.Range("E4").Value = 0
$ws.Range("J4").Value = 14.1116
$ws.Range("R4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AN4").Value = 11.1182
$ws.Range("AO4").Value = 0
$ws.Range("D5").Value = 14
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 11.2
$ws.Range("R5").Value = 2.4
$ws.Range("S5").Value = 0
$ws.Range("AC5").Value = 16.8
$ws.Range("AL5").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("P6").Value = 43.0917
$ws.Range("T6").Value = 35.5677
$ws.Range("Y6").Value = 0
$ws.Range("Z6").Value = 38.3037
$ws.Range("G7").Value = 0
$ws.Range("I7").Value = 53.8117
$ws.Range("Z7").Value = 0
$ws.Range("AG7").Value = 19.7309
$ws.Range("AM7").Value = 52.0179
$ws.Range("H8").Value = 0
$ws.Range("AA8").Value = 0.8038999999999999
$ws.Range("AB8").Value = 75.56270000000001
$ws.Range("AD8").Value = 0
$ws.Range("AI8").Value = 47.4277
$ws.Range("C9").Value = 0
$ws.Range("E9").Value = 38.4615
$ws.Range("P9").Value = 0
$ws.Range("Q9").Value = 70.8899
$ws.Range("U9").Value = 69.38160000000001
$ws.Range("X9").Value = 0
$ws.Range("AA9").Value = 0
$ws.Range("AI9").Value = 0
$ws.Range("F10").Value = 0.8811
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 80.17619999999999
$ws.Range("O10").Value = 0
$ws.Range("Q10").Value = 0
$ws.Range("Y10").Value = 24.6696
$ws.Range("C11").Value = 71.6216
$ws.Range("H11").Value = 8.1081
$ws.Range("L11").Value = 0
$ws.Range("U11").Value = 0
$ws.Range("V11").Value = 0
$ws.Range("AC11").Value = 0
$ws.Range("AD11").Value = 98.6486
$ws.Range("W12").Value = 0
$ws.Range("X12").Value = 35.9053
$ws.Range("AE12").Value = 42.7807
$ws.Range("AJ12").Value = 0
$ws.Range("AL12").Value = 33.6134
$ws.Range("AO12").Value = 61.8793
$ws.Range("K13").Value = 29.656
$ws.Range("M13").Value = 107.9478
$ws.Range("N13").Value = 0
$ws.Range("AE13").Value = 0
$ws.Range("AF13").Value = 27.2835
$ws.Range("AN13").Value = 0
$ws.Range("N14").Value = 39.604
$ws.Range("S14").Value = 53.9054
$ws.Range("V14").Value = 68.2068
$ws.Range("AK14").Value = 0
$ws.Range("AM14").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("O15").Value = 15.3846
$ws.Range("AG15").Value = 0
$ws.Range("AH15").Value = 1.2821
$ws.Range("AK15").Value = 94.87179999999999
